# Update "Training Dashboard" with the new progress-as-of date (04-Nov-2025).
# For every data row (3-30):
#   - column H (PERIOD TO EXPIRE) decreases by 1 day
#   - column I (LAST UPDATE) moves from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 30; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)   # column H
    $updateCell = $ws.Cells.Item($row, 9)   # column I

    $oldPeriod = $periodCell.Value()
    $periodCell.Value = $oldPeriod - 1

    # Force the cell to stay a plain text value (matching the existing
    # literal date string) instead of Excel auto-coercing it into a
    # real date serial number.
    $updateCell.NumberFormat = "@"
    $updateCell.Value = "04-Nov-2025"
}
